{"js": "// Rewrite the seven addition/subtraction word problems with new, more\n// interesting scenarios, and lengthen the answer blanks. The \"Number\n// sentence:\" prompt line is identical across every problem, so it is\n// replaced once (search finds every occurrence) for all paragraphs.\n\nasync function replaceAll(body, findText, replaceText) {\n  const results = body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    console.log(\"WARNING: text not found: \" + findText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// --- Problem 1 ---\nawait replaceAll(body,\n  \"Lily found eight shells on the beach. She gave three shells to her friend. How many shells does Lily have now?\",\n  \"Barnaby the bear found seven delicious honeycombs. He ate three of them. How many honeycombs does Barnaby have left?\");\nawait replaceAll(body,\n  \"Answer: Lily has ____________ shells now.\",\n  \"Answer: Barnaby has ________________________________________ honeycombs left.\");\n\n// --- Problem 2 ---\nawait replaceAll(body,\n  \"Ben has nine toy cars. He gets four more for his birthday. How many toy cars does Ben have in total?\",\n  \"Penelope had nine shiny buttons. She found four more under the sofa. How many buttons does Penelope have in total?\");\nawait replaceAll(body,\n  \"Answer: Ben has ____________ toy cars in total.\",\n  \"Answer: Penelope has ________________________________________ buttons in total.\");\n\n// --- Problem 3 ---\nawait replaceAll(body,\n  \"Sarah baked twelve cookies. Her brother ate one cookie. How many cookies are left?\",\n  \"Professor Bumble bought eight bouncing beans. He gave one to his friend. How many bouncing beans does Professor Bumble have now?\");\nawait replaceAll(body,\n  \"Answer: There are ____________ cookies left.\",\n  \"Answer: Professor Bumble has ________________________________________ bouncing beans now.\");\n\n// --- Problem 4 ---\nawait replaceAll(body,\n  \"David has seven pencils. Emily gives him five more. How many pencils does David have now?\",\n  \"Flora saw twelve fluffy sheep in a field. Two of them were eating dandelions. How many sheep were not eating dandelions?\");\nawait replaceAll(body,\n  \"Answer: David now has ____________ pencils.\",\n  \"Answer: ________________________________________ sheep were not eating dandelions.\");\n\n// --- Problem 5 ---\nawait replaceAll(body,\n  \"A farmer has ten cows in the field. He buys five more cows. How many cows does he have now?\",\n  \"Kevin the carrot had six orange crayons. He received three more crayons for his birthday. How many crayons does Kevin have?\");\nawait replaceAll(body,\n  \"Answer: The farmer now has ____________ cows.\",\n  \"Answer: Kevin has ________________________________________ crayons.\");\n\n// --- Problem 6 ---\nawait replaceAll(body,\n  \"There are fifteen children in a class. Three children are absent today. How many children are in the class today?\",\n  \"Brenda the badger baked fifteen cupcakes. She ate five of them. How many cupcakes are left?\");\nawait replaceAll(body,\n  \"Answer: There are ____________ children in the class today.\",\n  \"Answer: There are ________________________________________ cupcakes left.\");\n\n// --- Problem 7 ---\nawait replaceAll(body,\n  \"A baker made six cakes. He sold two of them. How many cakes does the baker have left?\",\n  \"Wilbur the worm collected four shiny pebbles. He found seven more. How many pebbles does Wilbur have altogether?\");\nawait replaceAll(body,\n  \"Answer: The baker has ____________ cakes left.\",\n  \"Answer: Wilbur has ________________________________________ pebbles altogether.\");\n\n// --- Shared \"Number sentence:\" blank line (identical across all 7 problems) ---\nawait replaceAll(body,\n  \"Number sentence: __________________________________________________\",\n  \"Number sentence: ______________________________________________________________________\");\n", "ps1": "# Rewrite the seven addition/subtraction word problems with new, more\n# interesting scenarios, and lengthen the answer blanks. The \"Number\n# sentence:\" prompt line is identical across every problem, so it is\n# replaced once for all paragraphs via wdReplaceAll.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $findText,      # FindText\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        Write-Output \"WARNING: text not found: $findText\"\n    }\n}\n\n# --- Problem 1 ---\nReplace-All \"Lily found eight shells on the beach. She gave three shells to her friend. How many shells does Lily have now?\" \"Barnaby the bear found seven delicious honeycombs. He ate three of them. How many honeycombs does Barnaby have left?\"\nReplace-All \"Answer: Lily has ____________ shells now.\" \"Answer: Barnaby has ________________________________________ honeycombs left.\"\n\n# --- Problem 2 ---\nReplace-All \"Ben has nine toy cars. He gets four more for his birthday. How many toy cars does Ben have in total?\" \"Penelope had nine shiny buttons. She found four more under the sofa. How many buttons does Penelope have in total?\"\nReplace-All \"Answer: Ben has ____________ toy cars in total.\" \"Answer: Penelope has ________________________________________ buttons in total.\"\n\n# --- Problem 3 ---\nReplace-All \"Sarah baked twelve cookies. Her brother ate one cookie. How many cookies are left?\" \"Professor Bumble bought eight bouncing beans. He gave one to his friend. How many bouncing beans does Professor Bumble have now?\"\nReplace-All \"Answer: There are ____________ cookies left.\" \"Answer: Professor Bumble has ________________________________________ bouncing beans now.\"\n\n# --- Problem 4 ---\nReplace-All \"David has seven pencils. Emily gives him five more. How many pencils does David have now?\" \"Flora saw twelve fluffy sheep in a field. Two of them were eating dandelions. How many sheep were not eating dandelions?\"\nReplace-All \"Answer: David now has ____________ pencils.\" \"Answer: ________________________________________ sheep were not eating dandelions.\"\n\n# --- Problem 5 ---\nReplace-All \"A farmer has ten cows in the field. He buys five more cows. How many cows does he have now?\" \"Kevin the carrot had six orange crayons. He received three more crayons for his birthday. How many crayons does Kevin have?\"\nReplace-All \"Answer: The farmer now has ____________ cows.\" \"Answer: Kevin has ________________________________________ crayons.\"\n\n# --- Problem 6 ---\nReplace-All \"There are fifteen children in a class. Three children are absent today. How many children are in the class today?\" \"Brenda the badger baked fifteen cupcakes. She ate five of them. How many cupcakes are left?\"\nReplace-All \"Answer: There are ____________ children in the class today.\" \"Answer: There are ________________________________________ cupcakes left.\"\n\n# --- Problem 7 ---\nReplace-All \"A baker made six cakes. He sold two of them. How many cakes does the baker have left?\" \"Wilbur the worm collected four shiny pebbles. He found seven more. How many pebbles does Wilbur have altogether?\"\nReplace-All \"Answer: The baker has ____________ cakes left.\" \"Answer: Wilbur has ________________________________________ pebbles altogether.\"\n\n# --- Shared \"Number sentence:\" blank line (identical across all 7 problems) ---\nReplace-All \"Number sentence: __________________________________________________\" \"Number sentence: ______________________________________________________________________\"\n"}
